$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E9").Value = 1

$ws.Range("E10").Select()
